$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the Fecha (D), Volumen (M), Precio minimo (N), Precio maximo (O),
# Precio promedio ponderado (P) and Precio $/Kg (S) values between rows 2 and 3
# so the row with the smaller volume/lower price (formerly row 3) becomes row 2
# and vice versa.

$columns = @("D", "M", "N", "O", "P", "S")

foreach ($col in $columns) {
    $cell2 = $ws.Range($col + "2")
    $cell3 = $ws.Range($col + "3")

    $value2 = $cell2.Value2
    $value3 = $cell3.Value2

    $cell2.Value = $value3
    $cell3.Value = $value2
}
